$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header row: "_old" suffix -> "_FV2310", "_new" suffix -> "_FV2404" ---
for ($c = 1; $c -le 10; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value = ($cell.Value2 -replace "_old$", "_FV2310")
}
for ($c = 12; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value = ($cell.Value2 -replace "_new$", "_FV2404")
}

# --- Turn the used range into an Excel Table (ListObject) with autofilter ---
$range = $ws.Range("A1:U74")
$tbl = $ws.ListObjects.Add(1, $range, $null, 1)
$tbl.Name = "Table1"

# --- Freeze the header row (pane split after row 1) ---
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
